# FEEDSTOCKS.xlsx edit: "Fixed thermal storage calculation and updated Conversion data"
#
# 1. OIL sheet: Opex_var_buy_USD2015kWh (column C, rows 2-25) 8.32E-2 -> 9.26E-2
# 2. NATURALGAS sheet: Opex_var_buy_USD2015kWh (column C, rows 2-25) 0.09 -> 0.218
# 3. Selection / active-sheet bookkeeping:
#      GRID        : selection G17 -> D19, no longer the active tab
#      OIL         : selection D6  -> D18
#      NATURALGAS  : selection C18 -> D10
#      SOLAR       : becomes the active tab (selection stays D47:E47)

$wb = $excel.ActiveWorkbook

# --- Update OIL Opex_var_buy_USD2015kWh values (column C) ---
$oil = $wb.Worksheets.Item("OIL")
$oil.Range("C2:C25").Value = 0.092600000000000002

# --- Update NATURALGAS Opex_var_buy_USD2015kWh values (column C) ---
$gas = $wb.Worksheets.Item("NATURALGAS")
$gas.Range("C2:C25").Value = 0.218

# --- Update selections on each sheet (this also makes that sheet the ---
# --- active tab while it runs, so SOLAR must be touched/activated last) ---
$grid = $wb.Worksheets.Item("GRID")
$grid.Range("D19").Select()

$oil.Range("D18").Select()

$gas.Range("D10").Select()

# SOLAR ends up as the final active sheet/tab, selection unchanged (D47:E47)
$solar = $wb.Worksheets.Item("SOLAR")
$solar.Activate()
